$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 4.447810488340769)
    3 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.21235251628842)
    4 = @(0.001754667048134761, 0.05231270169004087, 16.98373111632243, 0.4998867070740569, 1, 17.53768519213466)
    5 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    6 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    8 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 0, 14.40014219143469)
    9 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 1, 28.30127388105354)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
